# Updated cryptos list (price/volume refresh + two ranking swaps) per commit
# "Updated cryptos list on Fri May  5 03:29:09 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" (column D) values look numeric (e.g. "29.261.98", "1.0000") but are
# stored as literal text in the source sheet - force text format first so
# Excel does not silently reinterpret/round them as numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.261.98'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.900.32'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.98'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3911'
$ws.Range('E8').Value = '  -0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07878'
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9886'
$ws.Range('E10').Value = '  -1.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.03'
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.929.35'
$ws.Range('E12').Value = '  +3.29%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.082'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.749'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06988'
$ws.Range('E15').Value = '  -0.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.38'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001001'
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('E20').Value = '  -0.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.271.39'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.314'
$ws.Range('E22').Value = '  -1.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.07'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.092'
$ws.Range('E24').Value = '  +1.72%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.35'
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.46'
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.995'
$ws.Range('E27').Value = '  +2.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '118.62'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.918'
$ws.Range('E29').Value = '  -4.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09367'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9055'
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.286'
$ws.Range('E32').Value = '  -1.67%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.327'
$ws.Range('E33').Value = '  -1.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.220'
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.181'
$ws.Range('E35').Value = '  +1.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05788'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02090'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.000'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.771'
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5718'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1786'
$ws.Range('E41').Value = '  -1.31%  '
$ws.Range('E42').Value = '  -2.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.95'
$ws.Range('E43').Value = '  -0.91%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5348'
$ws.Range('E44').Value = '  -1.61%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.201'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07046'
$ws.Range('E46').Value = '  -0.86%  '
$ws.Range('B47').Value = 'MXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.601'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.859'
$ws.Range('E48').Value = '  -1.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '113.20'
$ws.Range('E49').Value = '  +0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.065'
$ws.Range('E50').Value = '  -2.84%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.27'
$ws.Range('E51').Value = '  -0.35%  '
